$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Action-Level2 values between row 2 and row 3 (test data sequence changed)
$n2 = $ws.Range("N2").Value
$n3 = $ws.Range("N3").Value
$ws.Range("N2").Value = $n3
$ws.Range("N3").Value = $n2

# Update view state to reflect the scrolled/selected cell from the diff
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("N5").Select()
